$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: Picture (图片) status changes from "In Progress" (进行中) to "Completed" (已完成)
$ws.Range("B9").Value = "已完成"
$ws.Range("B9").Interior.Color = 5287936

# Row 10: Popup Button (弹出按钮) status changes from "Not Started" (未开始) to "In Progress" (进行中)
$ws.Range("B10").Value = "进行中"
$ws.Range("B10").Interior.Color = 65535

# Update active cell selection to B9
$ws.Range("B9").Select()
